$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9: update title and link
$ws.Range("D9").Value = "기초 교육이 실패하면 벌어지는 일"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/when-basic-education-fails/#utm_source=rss&utm_medium=rss&utm_campaign=when-basic-education-fails"

# Row 27: update title only
$ws.Range("D27").Value = "TensorFlow Custom Op으로 데이터 변환 최적화하기"

# Row 44: update title and link
$ws.Range("D44").Value = "Meta-Learning과 MAML의 개념 정리"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/95"
